# Re-sort the worksheets so that "总计" (summary) comes first and
# "2021-Q1" (detail) comes second. This matches the commit:
# "update data with resort sheetname"

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")
$wsQ1 = $wb.Worksheets.Item("2021-Q1")

# Move "总计" to be before "2021-Q1", effectively making it the first sheet
$wsTotal.Move($wsQ1)

# Moving a sheet makes it active; restore "2021-Q1" as the selected/active
# sheet to match the original tab-selection state. Re-fetch the reference
# by name since the object handle captured before the move is stale.
$wsQ1 = $wb.Worksheets.Item("2021-Q1")
$wsQ1.Activate()
